$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New correlation matrix values (province sentiment correlation, pairplot update)
$values = @{
    "B2" = 1;                  "C2" = 0.3807625718760089;  "D2" = 0.3029493016797381;
    "E2" = 0.8249848017418434; "F2" = 0.5381002104292207;  "G2" = -0.5925727588665705;
    "H2" = 0.8407293522681448;

    "B3" = 0.3807625718760089; "C3" = 1;                   "D3" = 0.9323058924499702;
    "E3" = 0.7075058901348336; "F3" = 0.6315326640533218;  "G3" = 0.3596764627439427;
    "H3" = 0.6056155656471309;

    "B4" = 0.3029493016797381; "C4" = 0.9323058924499702;  "D4" = 1;
    "E4" = 0.6569309222873292; "F4" = 0.6043609692055382;  "G4" = 0.2820080977645044;
    "H4" = 0.5094350468276497;

    "B5" = 0.8249848017418434; "C5" = 0.7075058901348336;  "D5" = 0.6569309222873292;
    "E5" = 1;                  "F5" = 0.753748537200926;   "G5" = -0.3177175246249234;
    "H5" = 0.9049110834301681;

    "B6" = 0.5381002104292207; "C6" = 0.6315326640533218;  "D6" = 0.6043609692055382;
    "E6" = 0.753748537200926;  "F6" = 1;                   "G6" = 0.1143943837807325;
    "H6" = 0.723991459312084;

    "B7" = -0.5925727588665705;"C7" = 0.3596764627439427;  "D7" = 0.2820080977645044;
    "E7" = -0.3177175246249234;"F7" = 0.1143943837807325;  "G7" = 1;
    "H7" = -0.3014711968983573;

    "B8" = 0.8407293522681448; "C8" = 0.6056155656471309;  "D8" = 0.5094350468276497;
    "E8" = 0.9049110834301681; "F8" = 0.723991459312084;   "G8" = -0.3014711968983573;
    "H8" = 1;
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
